# Updates the cryptos price list (Price / Volume(1h) columns) to match
# the refreshed GitHub Actions data snapshot, including the ApeXProtocol /
# VeChain row swap at rows 44-45.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.486.29"
$ws.Range("E2").Value = "  -3.00%  "
$ws.Range("D3").Value = "2.218.22"
$ws.Range("E3").Value = "  -6.62%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'297.16"
$ws.Range("E5").Value = "  -4.64%  "
$ws.Range("D6").Value = "'82.77"
$ws.Range("E6").Value = "  -4.92%  "
$ws.Range("E7").Value = "  -3.92%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  -4.92%  "
$ws.Range("E10").Value = "  -7.81%  "
$ws.Range("D11").Value = "'29.18"
$ws.Range("E11").Value = "  -4.09%  "
$ws.Range("D12").Value = "'47.71"
$ws.Range("E12").Value = "  -9.62%  "
$ws.Range("E13").Value = "  -2.31%  "
$ws.Range("D14").Value = "2.566.02"
$ws.Range("E14").Value = "  -6.31%  "
$ws.Range("D15").Value = "'6.30"
$ws.Range("E15").Value = "  -3.61%  "
$ws.Range("D16").Value = "'14.10"
$ws.Range("E16").Value = "  -6.00%  "
$ws.Range("D17").Value = "2.214.37"
$ws.Range("E17").Value = "  -6.69%  "
$ws.Range("D18").Value = "'0.717"
$ws.Range("E18").Value = "  -5.57%  "
$ws.Range("D19").Value = "39.401.80"
$ws.Range("E19").Value = "  -3.08%  "
$ws.Range("E20").Value = "  -4.07%  "
$ws.Range("D21").Value = "'5.71"
$ws.Range("E21").Value = "  -7.05%  "
$ws.Range("D22").Value = "'65.10"
$ws.Range("E22").Value = "  -4.98%  "
$ws.Range("E23").Value = "  -4.23%  "
$ws.Range("D24").Value = "'227.81"
$ws.Range("E24").Value = "  -3.25%  "
$ws.Range("E25").Value = "  -0.18%  "
$ws.Range("D26").Value = "'2.41"
$ws.Range("E26").Value = "  -6.80%  "
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("D28").Value = "'22.63"
$ws.Range("E28").Value = "  -4.79%  "
$ws.Range("D29").Value = "'2.17"
$ws.Range("E29").Value = "  +0.65%  "
$ws.Range("D30").Value = "'9.09"
$ws.Range("E30").Value = "  -1.45%  "
$ws.Range("D31").Value = "'149.24"
$ws.Range("E31").Value = "  -3.19%  "
$ws.Range("D32").Value = "'31.86"
$ws.Range("E32").Value = "  -6.90%  "
$ws.Range("D33").Value = "'0.999"
$ws.Range("E33").Value = "  -0.23%  "
$ws.Range("D34").Value = "'4.85"
$ws.Range("E34").Value = "  -6.60%  "
$ws.Range("D35").Value = "'0.0695"
$ws.Range("E35").Value = "  -4.66%  "
$ws.Range("E36").Value = "  -3.29%  "
$ws.Range("E37").Value = "  -3.62%  "
$ws.Range("D38").Value = "'0.0972"
$ws.Range("E38").Value = "  -2.87%  "
$ws.Range("D39").Value = "'15.24"
$ws.Range("E39").Value = "  -4.17%  "
$ws.Range("D40").Value = "'2.63"
$ws.Range("E40").Value = "  -5.66%  "
$ws.Range("E41").Value = "  -3.58%  "
$ws.Range("D42").Value = "'3.63"
$ws.Range("E42").Value = "  -5.38%  "
$ws.Range("D43").Value = "1.907.46"
$ws.Range("E43").Value = "  -2.93%  "
$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D44").Value = "'2.04"
$ws.Range("E44").Value = "  -14.94%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "'0.0259"
$ws.Range("E45").Value = "  -3.84%  "
$ws.Range("D46").Value = "'9.03"
$ws.Range("E46").Value = "  -3.36%  "
$ws.Range("D47").Value = "'16.11"
$ws.Range("E47").Value = "  -9.18%  "
$ws.Range("D48").Value = "'2.62"
$ws.Range("E48").Value = "  -2.82%  "
$ws.Range("D49").Value = "2.435.16"
$ws.Range("E49").Value = "  -6.31%  "
$ws.Range("D50").Value = "'70.78"
$ws.Range("E50").Value = "  -1.48%  "
$ws.Range("D51").Value = "'87.27"
$ws.Range("E51").Value = "  -6.35%  "
